$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Update the "Splněno?" status column (E) for two tasks:
#  - E14 ("Rozsah prezentace ...") changes from "ne" to "WIP"
#  - E24 ("Databázi vyexportujte ... SQL") changes from "ne" to "ano"
$ws.Range("E14").Value = "WIP"
$ws.Range("E24").Value = "ano"

# Restore the selected cell to E14, matching the saved view state.
$ws.Range("E14").Select()

$wb.Save()
